# Auto-generated: bulk update of currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across multiple sheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1050
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -6224
$ws.Range("H15").Value = 144.75
$ws.Range("I15").Value = 144.75
$ws.Range("K15").Value = 434.25
$ws.Range("M15").Value = -265.25
$ws.Range("H80").Value = 5830.8696
$ws.Range("I80").Value = 1044
$ws.Range("J80").Value = 9513.076999999999
$ws.Range("K80").Value = 3132
$ws.Range("L80").Value = 28539.231
$ws.Range("M80").Value = -2134
$ws.Range("N80").Value = -30535.231
$ws.Range("H83").Value = 5830.8696
$ws.Range("I83").Value = 1044
$ws.Range("J83").Value = 9513.076999999999
$ws.Range("K83").Value = 9396
$ws.Range("L83").Value = 85617.693
$ws.Range("M83").Value = -4404
$ws.Range("N83").Value = -95601.693
$ws.Range("H88").Value = 1371.8948
$ws.Range("I88").Value = 1250
$ws.Range("J88").Value = 1460.5454
$ws.Range("K88").Value = 1250
$ws.Range("L88").Value = 1460.5454
$ws.Range("M88").Value = -844
$ws.Range("N88").Value = -2272.5454
$ws.Range("H91").Value = 1371.8948
$ws.Range("I91").Value = 1250
$ws.Range("J91").Value = 1460.5454
$ws.Range("K91").Value = 1250
$ws.Range("L91").Value = 1460.5454
$ws.Range("M91").Value = 154
$ws.Range("N91").Value = -4268.5454
$ws.Range("H112").Value = 1758.0454
$ws.Range("J112").Value = 1893.85
$ws.Range("L112").Value = 5681.549999999999
$ws.Range("N112").Value = -7897.549999999999
$ws.Range("H113").Value = 3073.45
$ws.Range("I113").Value = 2887.7144
$ws.Range("K113").Value = 2887.7144
$ws.Range("M113").Value = 366.2856000000002
$ws.Range("H138").Value = 1448263.9
$ws.Range("I138").Value = 2765.8572
$ws.Range("J138").Value = 1861263.4
$ws.Range("K138").Value = 8297.571599999999
$ws.Range("L138").Value = 5583790.199999999
$ws.Range("M138").Value = -3157.571599999999
$ws.Range("N138").Value = -5594070.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 33771.4
$ws.Range("J124").Value = 33771.4
$ws.Range("L124").Value = 33771.4
$ws.Range("N124").Value = -43591.4
$ws.Range("H125").Value = 67527
$ws.Range("J125").Value = 67527
$ws.Range("L125").Value = 67527
$ws.Range("N125").Value = -77367

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2642.913
$ws.Range("I86").Value = 2685.6667
$ws.Range("J86").Value = 2194
$ws.Range("K86").Value = 2685.6667
$ws.Range("L86").Value = 2194
$ws.Range("M86").Value = -1562.6667
$ws.Range("N86").Value = -4440
$ws.Range("H89").Value = 2642.913
$ws.Range("I89").Value = 2685.6667
$ws.Range("J89").Value = 2194
$ws.Range("K89").Value = 13428.3335
$ws.Range("L89").Value = 10970
$ws.Range("M89").Value = -7812.333500000001
$ws.Range("N89").Value = -22202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 518
$ws.Range("I105").Value = 518
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 518
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1229
$ws.Range("N105").Value = ""
$ws.Range("H132").Value = 2768.0435
$ws.Range("I132").Value = 2432.5293
$ws.Range("J132").Value = 3718.6667
$ws.Range("K132").Value = 7297.5879
$ws.Range("L132").Value = 11156.0001
$ws.Range("M132").Value = -4767.5879
$ws.Range("N132").Value = -16216.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("K7").Value = 600
$ws.Range("M7").Value = -488
$ws.Range("H33").Value = 555.9286
$ws.Range("I33").Value = 538.4
$ws.Range("J33").Value = 599.75
$ws.Range("K33").Value = 3230.4
$ws.Range("L33").Value = 3598.5
$ws.Range("M33").Value = -2947.4
$ws.Range("N33").Value = -4164.5
$ws.Range("H75").Value = 1161.5
$ws.Range("I75").Value = 613.1429000000001
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 1839.4287
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -841.4287000000002
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 1161.5
$ws.Range("I78").Value = 613.1429000000001
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 5518.2861
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -526.2861000000003
$ws.Range("N78").Value = -54984
$ws.Range("H107").Value = 2646611
$ws.Range("I107").Value = 13889336
$ws.Range("J107").Value = 1263.7354
$ws.Range("K107").Value = 41668008
$ws.Range("L107").Value = 3791.2062
$ws.Range("M107").Value = -41666088
$ws.Range("N107").Value = -7631.206200000001
$ws.Range("H129").Value = 2950.9092
$ws.Range("I129").Value = 4160
$ws.Range("K129").Value = 12480
$ws.Range("M129").Value = -7480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6775.5557
$ws.Range("I102").Value = 6200
$ws.Range("J102").Value = 7495
$ws.Range("K102").Value = 6200
$ws.Range("L102").Value = 7495
$ws.Range("M102").Value = -4578
$ws.Range("N102").Value = -10739
$ws.Range("H109").Value = 29285
$ws.Range("J109").Value = 29285
$ws.Range("L109").Value = 29285
$ws.Range("N109").Value = -31365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3823.3
$ws.Range("I40").Value = 3631.1667
$ws.Range("K40").Value = 3631.1667
$ws.Range("M40").Value = -3495.1667
$ws.Range("H133").Value = 55990
$ws.Range("J133").Value = 55990
$ws.Range("L133").Value = 55990
$ws.Range("N133").Value = -61050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2895.3684
$ws.Range("I62").Value = 2672.5
$ws.Range("J62").Value = 2998.2307
$ws.Range("K62").Value = 2672.5
$ws.Range("L62").Value = 2998.2307
$ws.Range("M62").Value = -2048.5
$ws.Range("N62").Value = -4246.2307
$ws.Range("H65").Value = 2895.3684
$ws.Range("I65").Value = 2672.5
$ws.Range("J65").Value = 2998.2307
$ws.Range("K65").Value = 13362.5
$ws.Range("L65").Value = 14991.1535
$ws.Range("M65").Value = -10242.5
$ws.Range("N65").Value = -21231.1535
$ws.Range("H136").Value = 5742.2266
$ws.Range("I136").Value = 3263.7273
$ws.Range("K136").Value = 9791.1819
$ws.Range("M136").Value = -7241.1819
